# Updated symbol list on Sat Dec 31 07:55:36 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price updates (Price column, D) ---
# A leading apostrophe is used so these numeric-looking strings are stored
# as text (matching the source data, which keeps trailing zeros etc.)
$ws.Range("D2").Value  = "'245.49"
$ws.Range("D3").Value  = "'25.52"
$ws.Range("D4").Value  = "'5.135"
$ws.Range("D6").Value  = "'6.484"
$ws.Range("D8").Value  = "'0.8170"
$ws.Range("D9").Value  = "'0.8427"

# --- Rows 11-19: coin ranking list shifted up by one position ---
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.02859"
$ws.Range("E11").Value = "10BitrueCoinBTR"

$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.09378"
$ws.Range("E12").Value = "11BitMartTokenBMX"

$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001515"
$ws.Range("E13").Value = "12BitForexTokenBF"

$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "'0.0005948"
$ws.Range("E14").Value = "13OneONE"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006146"
$ws.Range("E15").Value = "14TigerCashTCH"

$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.532"
$ws.Range("E16").Value = "15LEOLEO"

$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").Value = "'2.023"
$ws.Range("E17").Value = "16BTSETokenBTSE"

$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").Value = "'0.3179"
$ws.Range("E18").Value = "17BitpandaEcosystemTokenBEST"

$ws.Range("B19").Value = "MandalaExchangeToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D19").Value = "'0.06961"
$ws.Range("E19").Value = "18MandalaExchangeTokenMDX"

# --- Remaining simple price (and a couple of label) updates ---
$ws.Range("D20").Value = "'0.03218"
$ws.Range("D22").Value = "'3.740"
$ws.Range("D23").Value = "'0.04697"
$ws.Range("D26").Value = "'0.004608"

$ws.Range("D27").Value = "'0.00009698"
$ws.Range("E27").Value = "26NitroExNTX"

$ws.Range("D40").Value = "'0.03656"

$ws.Range("D41").Value = "'0.006148"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

$ws.Range("D44").Value = "'0.007772"
$ws.Range("D45").Value = "'0.00005314"

$ws.Range("D47").Value = "'0.1335"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

$ws.Range("D48").Value = "'0.002125"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D50").Value = "'0.0002000"
